$d = $word.ActiveDocument

$replacements = @(
    @("772÷6=", "221÷7="),
    @("416÷7=", "552÷5="),
    @("687÷7=", "285÷5="),
    @("105÷9=", "993÷9="),
    @("570÷9=", "798÷6="),
    @("212÷6=", "843÷2="),
    @("647÷2=", "355÷9="),
    @("560÷2=", "561÷3="),
    @("245÷3=", "647÷3="),
    @("930÷2=", "890÷5="),
    @("143÷3=", "292÷3="),
    @("735÷8=", "770÷3="),
    @("274÷9=", "283÷2="),
    @("991÷4=", "231÷9="),
    @("284÷3=", "837÷5="),
    @("222÷7=", "290÷2="),
    @("200÷6=", "807÷2="),
    @("596÷4=", "773÷3="),
    @("860÷7=", "690÷7="),
    @("971÷9=", "274÷4="),
    @("507÷4=", "476÷3="),
    @("705÷2=", "867÷2="),
    @("940÷6=", "268÷9="),
    @("133÷3=", "885÷4="),
    @("151÷2=", "684÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
